$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in a score that was missing for Будык Захар (row 8)
$ws.Range("C8").Value = 5

# Clear the computed "total" (N) and "grade" (P) helper columns for every
# student row - this removes the SUM formulas in N5:N29 and the literal
# values in P5:P29 while leaving the N6:N29 cells' formatting (style 6) in
# place, matching how a user would select the range and press Delete.
$ws.Range("N5:P29").ClearContents()

# Reproduce the new selection left behind by that delete operation.
$ws.Range("N4:P33").Select()

Write-Host "done"
